# Update: "Fruta / hortaliza, semanal"
# The data rows (2-11) had their Fecha (D), Volumen (J), Precio minimo (K),
# Precio maximo (L), Precio promedio ponderado (M) and Precio $/Kg (P)
# values reshuffled among the rows (a row permutation of the weekly
# price records). Capture the original values first, then write them
# back out in the new row order so every row ends up with the values
# that originally belonged to a different row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: target row -> source row (values that row should receive
# come from this original row number).
$rowMap = @{
    2  = 5
    3  = 11
    4  = 10
    5  = 4
    6  = 7
    7  = 9
    8  = 3
    9  = 6
    10 = 2
    11 = 8
}

$cols = @("D", "J", "K", "L", "M", "P")

# Snapshot the original values for the affected columns/rows before
# writing anything, since several rows are both sources and targets.
$orig = @{}
foreach ($r in 2..11) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $orig[$r] = $rowVals
}

foreach ($targetRow in $rowMap.Keys) {
    $sourceRow = $rowMap[$targetRow]
    $sourceVals = $orig[$sourceRow]
    foreach ($c in $cols) {
        $ws.Range("$c$targetRow").Value = $sourceVals[$c]
    }
}
